# Update RFMA recurrence/retention metrics for Dados_ADD_PF (metricas_recorrencia_anual)
# Columns: A=ano, B=ano_obj, C=total_customers, D=returning_customers, E=new_customers,
#          F=retention_rate, G=new_rate, H=returning_rate

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 2021
$ws.Range("C2").Value = 481
$ws.Range("D2").Value = 64
$ws.Range("E2").Value = 417
$ws.Range("F2").Value = 59.81308411214953
$ws.Range("G2").Value = 86.69438669438669
$ws.Range("H2").Value = 13.30561330561331

# Row 3 - 2022
$ws.Range("C3").Value = 359
$ws.Range("D3").Value = 164
$ws.Range("E3").Value = 195
$ws.Range("F3").Value = 34.0956340956341
$ws.Range("G3").Value = 54.31754874651811
$ws.Range("H3").Value = 45.68245125348189

# Row 4 - 2023
$ws.Range("C4").Value = 318
$ws.Range("D4").Value = 183
$ws.Range("E4").Value = 135
$ws.Range("F4").Value = 50.97493036211699
$ws.Range("G4").Value = 42.45283018867924
$ws.Range("H4").Value = 57.54716981132076

# Row 5 - 2024
$ws.Range("C5").Value = 449
$ws.Range("D5").Value = 228
$ws.Range("E5").Value = 221
$ws.Range("F5").Value = 71.69811320754717
$ws.Range("G5").Value = 49.22048997772828
$ws.Range("H5").Value = 50.77951002227172

# Row 6 - 2025
$ws.Range("C6").Value = 335
$ws.Range("D6").Value = 277
$ws.Range("E6").Value = 58
$ws.Range("F6").Value = 61.69265033407573
$ws.Range("G6").Value = 17.3134328358209
$ws.Range("H6").Value = 82.68656716417911
